# Actualización automática 2025-07-28 16:15:09
#
# Insert a new advisor row "DANIELA ELIZABETH BECERRA BECERRA" right above
# "EQUISAB S.A." (row 7) on both the "VENTAS POR GRUPO" and "VENTA MENSUAL"
# sheets. All existing rows from 7 downward shift down by one; the trailing
# summary/totals row (previously row 17) moves to row 18 and its "X de 15"
# labels are refreshed to "X de 16" to reflect the new row count.

$wb = $excel.ActiveWorkbook

$newName = "DANIELA ELIZABETH BECERRA BECERRA"
$office  = "OFICINA-CATAECSA"

# --- Sheet 1: "VENTAS POR GRUPO" (data columns C..R, summary row moves 17 -> 18)
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("A7").EntireRow.Insert()
$ws1.Range("A7").Value = $office
$ws1.Range("B7").Value = $newName
$ws1.Range("C7:R7").Value = 0

for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(18, $col)
    $cell.Value = $cell.Text.Replace("de 15", "de 16")
}

# --- Sheet 2: "VENTA MENSUAL" (data columns C..G, summary row moves 17 -> 18)
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("A7").EntireRow.Insert()
$ws2.Range("A7").Value = $office
$ws2.Range("B7").Value = $newName
$ws2.Range("C7:G7").Value = 0
